# Remove the paragraphs that hold the two inline pictures at the top of
# the document (the "27.jpg" / "26.jpg" drawings), leaving the trailing
# empty paragraph and the section properties untouched.
#
# Walk InlineShapes from the end backwards so removing one shape's
# paragraph never invalidates the indices of the ones still to be
# processed.

$d = $word.ActiveDocument

for ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {
    $shape = $d.InlineShapes.Item($i)
    $para = $shape.Range.Paragraphs.Item(1)
    $r = $d.Range($para.Range.Start, $para.Range.End)
    $r.Delete()
}
